$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New gyroscope readings (x, y, z) covering rows 2-31: the oldest two
# samples were dropped and twelve new samples appended (sliding window
# of streaming sensor data).
$numRows = 30
$numCols = 3
$data = New-Object 'object[,]' $numRows,$numCols
$data[0,0] = -0.2800817191600799
$data[0,1] = -1.379181027412415
$data[0,2] = 0.1039998084306716
$data[1,0] = 0.2446515262126922
$data[1,1] = 1.389260292053223
$data[1,2] = -0.0532979927957057
$data[2,0] = 1.946674823760986
$data[2,1] = 6.569244861602783
$data[2,2] = -1.187216639518738
$data[3,0] = 1.043357849121094
$data[3,1] = 2.923295736312866
$data[3,2] = 0.2732094824314117
$data[4,0] = 1.18660569190979
$data[4,1] = 1.866804242134094
$data[4,2] = 0.970054030418396
$data[5,0] = 0.1301143020391464
$data[5,1] = 2.222785949707031
$data[5,2] = -0.0280998013913631
$data[6,0] = -0.0589484944939613
$data[6,1] = -0.4569272100925445
$data[6,2] = -0.4526511430740356
$data[7,0] = 1.018770456314087
$data[7,1] = -5.163338661193848
$data[7,2] = -2.075872898101806
$data[8,0] = -1.822669148445129
$data[8,1] = -3.153133630752563
$data[8,2] = 0.7982481718063354
$data[9,0] = -0.737772524356842
$data[9,1] = -2.180025339126587
$data[9,2] = -0.1867720484733581
$data[10,0] = -0.2316706478595733
$data[10,1] = -2.678491353988647
$data[10,2] = -0.6579018831253052
$data[11,0] = -0.2446515262126922
$data[11,1] = 1.025031924247742
$data[11,2] = 0.4431827366352081
$data[12,0] = 0.9764680862426758
$data[12,1] = 5.428759574890137
$data[12,2] = 0.0936150997877121
$data[13,0] = 1.343140006065369
$data[13,1] = 3.488956928253174
$data[13,2] = 0.9002626538276672
$data[14,0] = -0.2557998299598694
$data[14,1] = 1.629635810852051
$data[14,2] = 1.214094638824463
$data[15,0] = -0.2167044430971145
$data[15,1] = 1.373072385787964
$data[15,2] = -0.3060434758663177
$data[16,0] = 0.1476766765117645
$data[16,1] = -0.5042692422866821
$data[16,2] = -0.403781920671463
$data[17,0] = 0.0123700210824608
$data[17,1] = -2.318233489990234
$data[17,2] = -0.6989825367927551
$data[18,0] = -1.0144944190979
$data[18,1] = -1.167210817337036
$data[18,2] = 0.6551529765129089
$data[19,0] = -0.3874412775039673
$data[19,1] = -0.7050912380218506
$data[19,2] = 0.0914770737290382
$data[20,0] = -0.1458440721035003
$data[20,1] = -0.3762930035591125
$data[20,2] = -0.0704022198915481
$data[21,0] = 0.2157881408929824
$data[21,1] = 0.3240640163421631
$data[21,2] = 0.0951422601938247
$data[22,0] = 0.0403171069920063
$data[22,1] = 0.1484402567148208
$data[22,2] = -0.0852157026529312
$data[23,0] = 0.1014036312699318
$data[23,1] = 0.3179553747177124
$data[23,2] = 0.0390953756868839
$data[24,0] = 0.052381694316864
$data[24,1] = 0.1099557429552078
$data[24,2] = 0.0681114718317985
$data[25,0] = 0.0522289797663688
$data[25,1] = -0.4196644127368927
$data[25,2] = 0.2273945808410644
$data[26,0] = 0.0937678143382072
$data[26,1] = -0.1565342247486114
$data[26,2] = 0.0675006061792373
$data[27,0] = -0.0591012127697467
$data[27,1] = 0.0331394411623477
$data[27,2] = 0.0291688162833452
$data[28,0] = -0.0047342055477201
$data[28,1] = 0.0655152946710586
$data[28,2] = -0.0163406450301408
$data[29,0] = -0.00534507073462
$data[29,1] = 0.0612392425537109
$data[29,2] = -0.0022907445672899

$startRow = 2
$range = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($startRow + $numRows - 1, $numCols))
$range.Value = $data
